$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: split the old combined "note - ..." text into a bold "NOTE - " label
# in column A, and move the remaining note text into column B.
$ws.Range("A4").Value = "NOTE - "
$ws.Range("A4").Font.Name = "Arial"
$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Font.Color = 0
$ws.Range("A4").Interior.Color = 16777215
$ws.Range("A4").HorizontalAlignment = -4131

$ws.Range("B4").Value = "there can be many options as much as you want, the last cell of each row should have the value for the answer."

# Row 5: new guidance about the required sheet name.
$ws.Range("B5").Value = "Sheet name should be ""Sheet1"""

# Row 6 (new): the "for more info" note, moved here and given the same
# white-fill / left-aligned styling as the new note rows (but not bold).
$ws.Range("B6").Value = "for more info, please view the sample."
$ws.Range("B6").Font.Name = "Arial"
$ws.Range("B6").Font.Bold = $false
$ws.Range("B6").Font.Color = 0
$ws.Range("B6").Interior.Color = 16777215
$ws.Range("B6").HorizontalAlignment = -4131

Write-Output "done"
